$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 67.5
$ws.Range("I53").Value = 67.5
$ws.Range("K53").Value = 67.5
$ws.Range("M53").Value = 569.5
$ws.Range("H82").Value = 8110.5
$ws.Range("I82").Value = 1222
$ws.Range("K82").Value = 3666
$ws.Range("M82").Value = -3260
$ws.Range("H85").Value = 8110.5
$ws.Range("I85").Value = 1222
$ws.Range("K85").Value = 3666
$ws.Range("M85").Value = -2262
$ws.Range("H107").Value = 559.5217
$ws.Range("I107").Value = 528.6
$ws.Range("K107").Value = 528.6
$ws.Range("M107").Value = 1391.4
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").ClearContents()
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = 0
$ws.Range("H131").Value = 834771.75
$ws.Range("I131").Value = 834771.75
$ws.Range("K131").Value = 2504315.25
$ws.Range("M131").Value = -2499275.25
$ws.Range("H132").Value = 10143.934
$ws.Range("I132").Value = 11630.77
$ws.Range("K132").Value = 34892.31
$ws.Range("M132").Value = -32362.31
$ws.Range("H137").Value = 4238.457
$ws.Range("I137").Value = 1421.7941
$ws.Range("K137").Value = 4265.3823
$ws.Range("M137").Value = -1715.3823
$ws.Range("H138").Value = 342494.94
$ws.Range("J138").Value = 443036.97
$ws.Range("L138").Value = 1329110.91
$ws.Range("N138").Value = -1339390.91

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6727.2666
$ws.Range("I32").Value = 6054.317
$ws.Range("K32").Value = 6054.317
$ws.Range("M32").Value = -5767.317
$ws.Range("H61").Value = 3540.5334
$ws.Range("I61").Value = 1551.7142
$ws.Range("K61").Value = 1551.7142
$ws.Range("M61").Value = -1339.7142
$ws.Range("H122").Value = 3034.5134
$ws.Range("I122").Value = 2819.697
$ws.Range("J122").Value = 4806.75
$ws.Range("K122").Value = 8459.091
$ws.Range("L122").Value = 14420.25
$ws.Range("M122").Value = -6009.091
$ws.Range("N122").Value = -19320.25
$ws.Range("H132").Value = 2465.318
$ws.Range("I132").Value = 2032.7273
$ws.Range("J132").Value = 3763.0908
$ws.Range("K132").Value = 6098.1819
$ws.Range("L132").Value = 11289.2724
$ws.Range("M132").Value = -3568.1819
$ws.Range("N132").Value = -16349.2724
$ws.Range("H136").Value = 3540.5334
$ws.Range("I136").Value = 1551.7142
$ws.Range("K136").Value = 4655.142599999999
$ws.Range("M136").Value = -2105.142599999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2733.0967
$ws.Range("I99").Value = 1927.25
$ws.Range("J99").Value = 5496
$ws.Range("K99").Value = 1927.25
$ws.Range("L99").Value = 5496
$ws.Range("M99").Value = -429.25
$ws.Range("N99").Value = -8492
$ws.Range("H107").Value = 1385.2565
$ws.Range("I107").Value = 1246.8928
$ws.Range("J107").Value = 1737.4546
$ws.Range("K107").Value = 1246.8928
$ws.Range("L107").Value = 1737.4546
$ws.Range("M107").Value = 673.1071999999999
$ws.Range("N107").Value = -5577.4546

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3556.7083
$ws.Range("I31").Value = 2864.3262
$ws.Range("J31").Value = 4781.6924
$ws.Range("K31").Value = 2864.3262
$ws.Range("L31").Value = 4781.6924
$ws.Range("M31").Value = -2569.3262
$ws.Range("N31").Value = -5371.6924
$ws.Range("H34").Value = 3556.7083
$ws.Range("I34").Value = 2864.3262
$ws.Range("J34").Value = 4781.6924
$ws.Range("K34").Value = 2864.3262
$ws.Range("L34").Value = 4781.6924
$ws.Range("M34").Value = -2662.3262
$ws.Range("N34").Value = -5185.6924
$ws.Range("H58").Value = 2317.4243
$ws.Range("I58").Value = 1333.0869
$ws.Range("K58").Value = 1333.0869
$ws.Range("M58").Value = -1130.0869
$ws.Range("H99").Value = 8468.4
$ws.Range("I99").Value = 9160.625
$ws.Range("J99").Value = 5699.5
$ws.Range("K99").Value = 9160.625
$ws.Range("L99").Value = 5699.5
$ws.Range("M99").Value = -7662.625
$ws.Range("N99").Value = -8695.5
$ws.Range("H126").Value = 8468.4
$ws.Range("I126").Value = 9160.625
$ws.Range("J126").Value = 5699.5
$ws.Range("K126").Value = 27481.875
$ws.Range("L126").Value = 17098.5
$ws.Range("M126").Value = -25011.875
$ws.Range("N126").Value = -22038.5
$ws.Range("H134").Value = 2378.0334
$ws.Range("I134").Value = 2086.7036
$ws.Range("K134").Value = 6260.110799999999
$ws.Range("M134").Value = -3725.110799999999
$ws.Range("H136").Value = 2317.4243
$ws.Range("I136").Value = 1333.0869
$ws.Range("K136").Value = 3999.2607
$ws.Range("M136").Value = -1449.2607

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 2919.3333
$ws.Range("J63").Value = 5374.5
$ws.Range("L63").Value = 16123.5
$ws.Range("N63").Value = -17621.5
$ws.Range("H66").Value = 2919.3333
$ws.Range("J66").Value = 5374.5
$ws.Range("L66").Value = 48370.5
$ws.Range("N66").Value = -55858.5
$ws.Range("H86").Value = 316.66666
$ws.Range("I86").Value = 280
$ws.Range("K86").Value = 840
$ws.Range("M86").Value = 346
$ws.Range("H89").Value = 316.66666
$ws.Range("I89").Value = 280
$ws.Range("K89").Value = 2520
$ws.Range("M89").Value = 3408
$ws.Range("H92").Value = 498
$ws.Range("I92").Value = 547.5
$ws.Range("J92").Value = 448.5
$ws.Range("K92").Value = 1642.5
$ws.Range("L92").Value = 1345.5
$ws.Range("M92").Value = -394.5
$ws.Range("N92").Value = -3841.5
$ws.Range("H117").Value = 2882.25
$ws.Range("I117").Value = 2764.5
$ws.Range("J117").Value = 3000
$ws.Range("K117").Value = 8293.5
$ws.Range("L117").Value = 9000
$ws.Range("M117").Value = -4851.5
$ws.Range("N117").Value = -15884

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1863.7084
$ws.Range("I102").Value = 1569.2632
$ws.Range("J102").Value = 2982.6
$ws.Range("K102").Value = 1569.2632
$ws.Range("L102").Value = 2982.6
$ws.Range("M102").Value = 52.7367999999999
$ws.Range("N102").Value = -6226.6
$ws.Range("H126").Value = 9791.666999999999
$ws.Range("I126").Value = 3959
$ws.Range("J126").Value = 21457
$ws.Range("K126").Value = 11877
$ws.Range("L126").Value = 64371
$ws.Range("M126").Value = -9407
$ws.Range("N126").Value = -69311
$ws.Range("H132").Value = 3851.1428
$ws.Range("I132").Value = 3167.6155
$ws.Range("J132").Value = 4961.875
$ws.Range("K132").Value = 9502.8465
$ws.Range("L132").Value = 14885.625
$ws.Range("M132").Value = -6972.8465
$ws.Range("N132").Value = -19945.625

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1567.9744
$ws.Range("I46").Value = 1993.75
$ws.Range("J46").Value = 1458.0968
$ws.Range("K46").Value = 1993.75
$ws.Range("L46").Value = 1458.0968
$ws.Range("M46").Value = -1805.75
$ws.Range("N46").Value = -1834.0968
$ws.Range("H122").Value = 3952.6785
$ws.Range("I122").Value = 3785.25
$ws.Range("J122").Value = 4371.25
$ws.Range("K122").Value = 11355.75
$ws.Range("L122").Value = 13113.75
$ws.Range("M122").Value = -8905.75
$ws.Range("N122").Value = -18013.75
$ws.Range("H132").Value = 4457.6177
$ws.Range("I132").Value = 2728.4
$ws.Range("J132").Value = 6927.9287
$ws.Range("K132").Value = 8185.200000000001
$ws.Range("L132").Value = 20783.7861
$ws.Range("M132").Value = -5655.200000000001
$ws.Range("N132").Value = -25843.7861
$ws.Range("H133").Value = 103995
$ws.Range("J133").Value = 103995
$ws.Range("L133").Value = 103995
$ws.Range("N133").Value = -109055

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1484.7858
$ws.Range("I126").Value = 1484.7858
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 4454.357400000001
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -1984.357400000001
$ws.Range("H132").Value = 3859.532
$ws.Range("I132").Value = 4030.611
$ws.Range("J132").Value = 3299.6365
$ws.Range("K132").Value = 12091.833
$ws.Range("L132").Value = 9898.9095
$ws.Range("M132").Value = -9561.832999999999
$ws.Range("N132").Value = -14958.9095
$ws.Range("H136").Value = 27780376
$ws.Range("I136").Value = 30304228
$ws.Range("K136").Value = 90912684
$ws.Range("M136").Value = -90910134
